$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header row (B1:E1): plain (non-numeric-looking) text -> just assign.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "Economic Performance"
$ws.Range("C1").Value = "Government Efficiency"
$ws.Range("D1").Value = "Business Efficiency"
$ws.Range("E1").Value = "Infrastructure"

# ------------------------------------------------------------------
# Existing rows 2-5: column A switches from numeric row-index to a
# text label ("score_2024"/"2020"/"2021"/"2022"), and B:E switch from
# text labels/positions to numeric-looking score strings. Every one of
# these must stay TEXT (matching the source inlineStr cells), so for
# any numeric-looking string we force NumberFormat "@" on that single
# cell *before* assigning the value (this keeps each column's existing
# style - bold/border on A, plain on B:E - and only adds the text
# format, reusing one shared style per column instead of creating a
# style per cell).
# ------------------------------------------------------------------

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "score_2024"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "75.0"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "49.8"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.1"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "73.7"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2020"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "26"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2021"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "5"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "28"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "10"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "6"

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2022"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "3"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "27"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "12"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "7"

# ------------------------------------------------------------------
# New rows 6-7: these cells don't exist yet, so first clone the
# formatting (bold/border/centered) of the existing column-A header
# style onto the new A-cell via a formats-only paste (this reuses the
# existing style id rather than minting a new one), then force the
# text number format + value on column A, and finally fill in the
# (plainly-styled) B:E cells with forced text format + value.
# ------------------------------------------------------------------

# Row 6
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2023"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "25"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "14"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6"

# Row 7
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Value = "2024"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "1"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "34"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "19"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "7"
